# Refresh the coin price / 1h-volume table with the latest scrape.
# (D-column price strings are quote-prefixed with a leading apostrophe so
#  Excel keeps numeric-looking text like "1.002" as text instead of coercing
#  it to a Number cell, matching the source data which is text.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'26.906.87"
$ws.Cells.Item(2, 5).Value = "  -0.95%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "'1.736.94"
$ws.Cells.Item(3, 5).Value = "  +0.99%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "'1.002"
$ws.Cells.Item(4, 5).Value = "  -0.04%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'310.99"
$ws.Cells.Item(5, 5).Value = "  -0.44%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  -0.06%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.4991"
$ws.Cells.Item(7, 5).Value = "  +8.56%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.3555"

# Row 9
$ws.Cells.Item(9, 4).Value = "'42.06"
$ws.Cells.Item(9, 5).Value = "  -1.23%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'0.07240"
$ws.Cells.Item(10, 5).Value = "  -0.25%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'1.058"
$ws.Cells.Item(11, 5).Value = "  +1.60%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -0.11%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'20.19"
$ws.Cells.Item(13, 5).Value = "  +2.21%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'5.934"
$ws.Cells.Item(14, 5).Value = "  +1.78%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'1.741.84"
$ws.Cells.Item(15, 5).Value = "  +1.08%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'6.827"
$ws.Cells.Item(16, 5).Value = "  -0.54%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'86.57"
$ws.Cells.Item(17, 5).Value = "  -3.15%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'0.00001033"
$ws.Cells.Item(18, 5).Value = "  -0.43%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'0.06395"
$ws.Cells.Item(19, 5).Value = "  +1.11%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -0.01%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'16.49"
$ws.Cells.Item(21, 5).Value = "  +0.26%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'5.707"
$ws.Cells.Item(22, 5).Value = "  +1.73%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'26.988.52"
$ws.Cells.Item(23, 5).Value = "  -0.85%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'11.32"
$ws.Cells.Item(24, 5).Value = "  +4.23%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'2.041"
$ws.Cells.Item(25, 5).Value = "  -4.26%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'154.32"
$ws.Cells.Item(26, 5).Value = "  -0.17%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'19.81"
$ws.Cells.Item(27, 5).Value = "  +3.00%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'1.938.52"
$ws.Cells.Item(28, 5).Value = "  +0.92%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'2.202"
$ws.Cells.Item(29, 5).Value = "  +2.56%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'119.58"
$ws.Cells.Item(30, 5).Value = "  +0.41%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'1.040"
$ws.Cells.Item(31, 5).Value = "  +1.53%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'0.09473"
$ws.Cells.Item(32, 5).Value = "  +4.28%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -0.15%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'5.348"
$ws.Cells.Item(34, 5).Value = "  +0.29%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'0.02188"
$ws.Cells.Item(35, 5).Value = "  -0.56%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'0.05842"
$ws.Cells.Item(36, 5).Value = "  -0.07%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'11.04"
$ws.Cells.Item(37, 5).Value = "  -0.08%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'1.426"
$ws.Cells.Item(38, 5).Value = "  +0.42%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'0.1995"
$ws.Cells.Item(39, 5).Value = "  +0.04%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'4.765"
$ws.Cells.Item(40, 5).Value = "  +0.80%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.6028"
$ws.Cells.Item(41, 5).Value = "  +1.62%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'1.105"
$ws.Cells.Item(42, 5).Value = "  -1.85%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "'7.608"
$ws.Cells.Item(43, 5).Value = "  +1.96%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +0.27%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'3.591"
$ws.Cells.Item(45, 5).Value = "  -0.04%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'0.5645"
$ws.Cells.Item(46, 5).Value = "  +0.20%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'120.05"
$ws.Cells.Item(47, 5).Value = "  +0.26%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'1.848"
$ws.Cells.Item(48, 5).Value = "  -0.87%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(49, 4).Value = "'0.06666"
$ws.Cells.Item(49, 5).Value = "  +0.11%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "EOS"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Cells.Item(50, 4).Value = "'1.101"
$ws.Cells.Item(50, 5).Value = "  +1.86%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  -0.01%  "
